$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 6, pushing the existing rows 6-12 down to 7-13
$ws.Rows("6:6").Insert()

# Populate the newly inserted row 6 with the new weekly record
$ws.Range("A6").Value = 11
$ws.Range("B6").Value = "Vega Monumental Concepción"
$ws.Range("C6").Value = "Bíobío"
$ws.Range("D6").Value = "2023-01-31"
$ws.Range("E6").Value = 8
$ws.Range("F6").Value = 100112052
$ws.Range("G6").Value = "Albahaca"
$ws.Range("H6").Value = "Sin especificar"
$ws.Range("I6").Value = "Primera"
$ws.Range("J6").Value = 70
$ws.Range("K6").Value = 1500
$ws.Range("L6").Value = 2000
$ws.Range("M6").Value = 1857
$ws.Range("N6").Value = "`$/docena de matas"
$ws.Range("O6").Value = "Región Metropolitana"
$ws.Range("P6").Value = 310
$ws.Range("Q6").Value = 6
$ws.Range("R6").Value = "Hortaliza"
